# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy the formatting from G1 (bold, bordered, centered header style)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# "Save" indicator (1/0) for each row, 2 through 52.
$saveValues = @(0,0,0,0,0,1,0,1,0,0,0,0,1,0,0,0,0,1,1,0,0,1,1,0,0,0,0,0,0,1,0,0,1,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,1,0,0)

$row = 2
foreach ($val in $saveValues) {
    $ws.Cells.Item($row, 8).Value = $val
    $row = $row + 1
}
